# "added functionality to pick best lot"
# Relabel the Mon/Tue/Wed/Thu/Fri time-slot column headers (row 1, columns
# D:W) so the leading zero before the hour digit is replaced with a dash,
# e.g. "Mon08" -> "Mon-8", "Mon010" -> "Mon-10", "Fri02" -> "Fri-2".
#
# All twenty header cells have to be rewritten together (not just the ones
# that visually change) so that every old "<Day><0><Hour>" shared string
# becomes fully unreferenced and the workbook naturally reuses the same
# shared-string slots for the new "<Day>-<Hour>" text instead of appending
# new entries at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Mon-8"
$ws.Range("E1").Value = "Mon-10"
$ws.Range("F1").Value = "Mon-12"
$ws.Range("G1").Value = "Mon-2"

$ws.Range("H1").Value = "Tue-8"
$ws.Range("I1").Value = "Tue-10"
$ws.Range("J1").Value = "Tue-12"
$ws.Range("K1").Value = "Tue-2"

$ws.Range("L1").Value = "Wed-8"
$ws.Range("M1").Value = "Wed-10"
$ws.Range("N1").Value = "Wed-12"
$ws.Range("O1").Value = "Wed-2"

$ws.Range("P1").Value = "Thu-8"
$ws.Range("Q1").Value = "Thu-10"
$ws.Range("R1").Value = "Thu-12"
$ws.Range("S1").Value = "Thu-2"

$ws.Range("T1").Value = "Fri-8"
$ws.Range("U1").Value = "Fri-10"
$ws.Range("V1").Value = "Fri-12"
$ws.Range("W1").Value = "Fri-2"

# Move the live selection from J12 to Q8, matching where the author's
# cursor ended up while picking the best lot.
$ws.Activate()
$ws.Range("Q8").Select()
